$wb = $excel.ActiveWorkbook

# ALC!row98: The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2668.0425
$ws.Range("I98").Value = 2330.8718
$ws.Range("J98").Value = 4311.75
$ws.Range("K98").Value = 2330.8718
$ws.Range("L98").Value = 4311.75
$ws.Range("M98").Value = -832.8717999999999
$ws.Range("N98").Value = -7307.75

# ALC!row122: Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2668.0425
$ws.Range("I122").Value = 2330.8718
$ws.Range("J122").Value = 4311.75
$ws.Range("K122").Value = 6992.6154
$ws.Range("L122").Value = 12935.25
$ws.Range("M122").Value = -4542.6154
$ws.Range("N122").Value = -17835.25

# ALC!row132: Fast-forwarding Flora / Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4044.0715
$ws.Range("I132").Value = 1329.5
$ws.Range("J132").Value = 20331.5
$ws.Range("K132").Value = 3988.5
$ws.Range("L132").Value = 60994.5
$ws.Range("M132").Value = -1458.5
$ws.Range("N132").Value = -66054.5

# ALC!row137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 64060.75
$ws.Range("I137").Value = 1757.1111
$ws.Range("J137").Value = 144165.42
$ws.Range("K137").Value = 5271.3333
$ws.Range("L137").Value = 432496.26
$ws.Range("M137").Value = -2721.3333
$ws.Range("N137").Value = -437596.26

# ARM!row16: Greavous Losses / Bronze Sabatons
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 686.3333
$ws.Range("I16").Value = 223.6
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 223.6
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 63.40000000000001
$ws.Range("N16").Value = -3574

# ARM!row23: A Well-rounded Crew / Iron Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 14999.5
$ws.Range("J23").Value = 14999.5
$ws.Range("L23").Value = 14999.5
$ws.Range("N23").Value = -15517.5

# ARM!row34: Insistent Sallets / Steel Sallet
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 30123
$ws.Range("J34").Value = 30123
$ws.Range("L34").Value = 30123
$ws.Range("N34").Value = -30665

# ARM!row45: Hollow Hallmarks / Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15396.5
$ws.Range("I45").Value = 15663.066
$ws.Range("J45").Value = 14952.223
$ws.Range("K45").Value = 15663.066
$ws.Range("L45").Value = 14952.223
$ws.Range("M45").Value = -15286.066
$ws.Range("N45").Value = -15706.223

# ARM!row110: Scheduled Maintenance / Deepgold Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1984.5883
$ws.Range("I110").Value = 1946.125
$ws.Range("K110").Value = 1946.125
$ws.Range("M110").Value = 98.875

# ARM!row122: Haste for High Durium / High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 32328.285
$ws.Range("I122").Value = 36883.168
$ws.Range("K122").Value = 110649.504
$ws.Range("M122").Value = -108199.504

# BSM!row20: Smelt and Dealt / Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12627.667
$ws.Range("I20").Value = 10235.111
$ws.Range("K20").Value = 10235.111
$ws.Range("M20").Value = -9988.111000000001

# BSM!row94: High Steal / High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1029.8334
$ws.Range("I94").Value = 815.86664
$ws.Range("J94").Value = 1386.4445
$ws.Range("K94").Value = 815.86664
$ws.Range("L94").Value = 1386.4445
$ws.Range("M94").Value = -364.86664
$ws.Range("N94").Value = -2288.4445

# CRP!row58: You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6459.05
$ws.Range("I58").Value = 6621.222
$ws.Range("K58").Value = 6621.222
$ws.Range("M58").Value = -6418.222

# CRP!row99: O Pine / Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1695.6
$ws.Range("I99").Value = 1394.5
$ws.Range("K99").Value = 1394.5
$ws.Range("M99").Value = 103.5

# CRP!row103: Spare a Rod and Spoil the Fishers / Gazelle Horn Fishing Rod
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 21868.6
$ws.Range("I103").Value = 22085.75
$ws.Range("K103").Value = 22085.75
$ws.Range("M103").Value = -20913.75

# CRP!row126: A Better Conductor / Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1695.6
$ws.Range("I126").Value = 1394.5
$ws.Range("K126").Value = 4183.5
$ws.Range("M126").Value = -1713.5

# CRP!row134: Wood You Be Quiet / Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2592.3235
$ws.Range("I134").Value = 2403.8667
$ws.Range("K134").Value = 7211.6001
$ws.Range("M134").Value = -4676.6001

# CRP!row136: Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6459.05
$ws.Range("I136").Value = 6621.222
$ws.Range("K136").Value = 19863.666
$ws.Range("M136").Value = -17313.666

# CUL!row37: I Love Lamprey / Eel Pie
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 500047500
$ws.Range("J37").Value = 500047500
$ws.Range("L37").Value = 1500142500
$ws.Range("N37").Value = -1500142724

# CUL!row115: Mixology / Blood Tomato Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 4031.4
$ws.Range("J115").Value = 9000
$ws.Range("L115").Value = 27000
$ws.Range("N115").Value = -29350

# CUL!row124: Bobbing for Compliments / Island Miq'abob
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 9228.6
$ws.Range("I124").Value = 849.5
$ws.Range("J124").Value = 14814.667
$ws.Range("K124").Value = 2548.5
$ws.Range("L124").Value = 44444.001
$ws.Range("M124").Value = 2361.5
$ws.Range("N124").Value = -54264.001

# CUL!row141: Ocean Explosion / Acqua Pazza
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3643
$ws.Range("I141").Value = 2128.7334
$ws.Range("K141").Value = 6386.2002
$ws.Range("M141").Value = -1206.2002

# GSM!row34: All Booked Up / Silver Magnifiers
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 20086.5
$ws.Range("I34").Value = 15000
$ws.Range("K34").Value = 15000
$ws.Range("M34").Value = -14732

# GSM!row76: The Monuments Mages / Hardsilver Magnifiers of Casting
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 20086.5
$ws.Range("I76").Value = 15000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14685

# GSM!row79: Deal with It (L) / Hardsilver Magnifiers of Casting
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 20086.5
$ws.Range("I79").Value = 15000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13908

# GSM!row102: Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1261.8462
$ws.Range("I102").Value = 767.1111
$ws.Range("K102").Value = 767.1111
$ws.Range("M102").Value = 854.8889

# LTW!row29: Hands On / Fingerless Goatskin Gloves
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10590

# LTW!row43: Subordinate Clause / Goatskin Choker
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 19978940
$ws.Range("I43").Value = 14023120
$ws.Range("K43").Value = 14023120
$ws.Range("M43").Value = -14022927

# LTW!row46: Supply Side Logic / Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 41074.637
$ws.Range("I46").Value = 72353.164
$ws.Range("J46").Value = 3540.4
$ws.Range("K46").Value = 72353.164
$ws.Range("L46").Value = 3540.4
$ws.Range("M46").Value = -72165.164
$ws.Range("N46").Value = -3916.4

# LTW!row61: Spelling Me Softly / Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11509448
$ws.Range("I61").Value = 12824372
$ws.Range("K61").Value = 12824372
$ws.Range("M61").Value = -12824170

# LTW!row113: Peace in Rest / Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 11509448
$ws.Range("I113").Value = 12824372
$ws.Range("K113").Value = 12824372
$ws.Range("M113").Value = -12822202

# LTW!row130: Generous Soles / Ophiotauroskin Boots of Healing
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 110000
$ws.Range("J130").Value = 110000
$ws.Range("L130").Value = 110000
$ws.Range("N130").Value = -120040

# LTW!row132: Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7909.3887
$ws.Range("I132").Value = 7904.353
$ws.Range("K132").Value = 23713.059
$ws.Range("M132").Value = -21183.059
